# Ajustando itens do template
#
# 1) "Convênio ou Contrato nº ..." -> "Convênio nº ..."
# 2) Collapse the split "{ind2}" / "{ativ2}" / "{totalCH2}" / "{#at2}{Dias}" /
#    "at2" template placeholders (each was typed across several runs) back
#    into single runs, without altering the visible text.
# 3) "Data:27/{mes}/2021" -> "Data: {final}"

$d = $word.ActiveDocument

# --- 1. Drop " ou Contrato" from the "Convênio ou Contrato nº ..." line ---
$r = $d.Content
$r.Find.Execute("Convênio ou Contrato", $true, $false, $false, $false, $false, `
                 $true, 1, $false, "Convênio", 2)

# --- 2. Re-merge the runs that spell out the template placeholders ---
# Doing a Find/Replace of the text with itself makes the engine coalesce the
# (identically-formatted) runs spanned by the match into a single run.
$r = $d.Content
$r.Find.Execute("{ind2}", $true, $false, $false, $false, $false, `
                 $true, 1, $false, "{ind2}", 2)

$r = $d.Content
$r.Find.Execute("{ativ2}", $true, $false, $false, $false, $false, `
                 $true, 1, $false, "{ativ2}", 2)

$r = $d.Content
$r.Find.Execute("{totalCH2}", $true, $false, $false, $false, $false, `
                 $true, 1, $false, "{totalCH2}", 2)

# This single pass fixes both "{#at" + "2" + "}{Dias}" -> "{#at2}{Dias}" and
# "at" + "2" -> "at2" inside "{CH}{/at2}" (the latter keeps its underline
# formatting separate from the neighbouring "{CH}{/" and "}" runs).
$r = $d.Content
$r.Find.Execute("at2", $true, $false, $false, $false, $false, `
                 $true, 1, $false, "at2", 2)

# --- 3. "Data:27/{mes}/2021" -> "Data: {final}" ---
$p = $d.Paragraphs.Last

# ":" becomes ": " (keeps "Data" as its own run, merges the added space into
# the existing ":" run).
$rngColon = $p.Range
$rngColon.Find.Execute(":", $true, $false, $false, $false, $false, `
                        $true, 1, $false, "")
$rngColon.Text = ": "

# Remove "27/{mes}/2021" and type the replacement as three separate runs
# ("{", "final", "}"), matching how the template placeholder is authored
# elsewhere in the document.
$rngTail = $p.Range
$rngTail.Find.Execute("27/{mes}/2021", $true, $false, $false, $false, $false, `
                       $true, 1, $false, "")
$rngTail.Text = ""
$rngTail.Collapse(0)
$rngTail.InsertAfter("{")
$rngTail.Collapse(0)
$rngTail.InsertAfter("final")
$rngTail.Collapse(0)
$rngTail.InsertAfter("}")
